# "Concluido o teste do metodo salvarCliente."
#
# The "SalvarCliente" test sheet had its expected-result column (H) filled
# in with placeholder/unrelated strings while the test was still in
# progress. Finishing the test run means the method under test now
# actually redirects to the real JSP views, so column H is updated to
# reflect the final, correct expected results:
#   - rows 3-8  (failure cases)  -> cliente/falha-cadastrar-cliente
#   - rows 9-14 (success cases)  -> cliente/sucesso-cadastro-cliente
#
# The workbook is also left with the "SalvarCliente" sheet active/selected
# (this was the sheet being worked on), with the cursor resting on C18.

$wb = $excel.ActiveWorkbook

$wsSalvarCliente = $wb.Worksheets.Item("SalvarCliente")

# Finish filling in the expected-result column now that the method's
# behaviour has been verified end-to-end.
$wsSalvarCliente.Range("H3:H8").Value  = "cliente/falha-cadastrar-cliente"
$wsSalvarCliente.Range("H9:H14").Value = "cliente/sucesso-cadastro-cliente"

# Leave the SalvarCliente sheet as the active/selected one, cursor on C18,
# since that's the test that was just completed.
$wsSalvarCliente.Activate()
$wsSalvarCliente.Range("C18").Select()
